# Apply 2024-10-18 data update to violent-crime-ytd workbook
# Updates column K (year 2024) values across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets, per the commit 'Add data for 2024-10-18'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 6363  # was 6339
$ws.Cells.Item(3, 11).Value = 6569  # was 6545
$ws.Cells.Item(4, 11).Value = 1367  # was 1358
$ws.Cells.Item(6, 11).Value = 7234  # was 7209
$ws.Cells.Item(7, 11).Value = 21997  # was 21915

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 396  # was 395
$ws.Cells.Item(3, 11).Value = 441  # was 436
$ws.Cells.Item(6, 11).Value = 487  # was 482
$ws.Cells.Item(7, 11).Value = 1447  # was 1436

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(4, 11).Value = 50  # was 49
$ws.Cells.Item(6, 11).Value = 296  # was 295
$ws.Cells.Item(7, 11).Value = 961  # was 959

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 123  # was 121
$ws.Cells.Item(7, 11).Value = 362  # was 360

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 211  # was 210
$ws.Cells.Item(3, 11).Value = 248  # was 245
$ws.Cells.Item(6, 11).Value = 223  # was 222
$ws.Cells.Item(7, 11).Value = 748  # was 743

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(6, 11).Value = 92  # was 91
$ws.Cells.Item(7, 11).Value = 362  # was 361

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(3, 11).Value = 26  # was 24
$ws.Cells.Item(7, 11).Value = 83  # was 81

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 191  # was 190
$ws.Cells.Item(5, 11).Value = 57  # was 56
$ws.Cells.Item(6, 11).Value = 157  # was 155
$ws.Cells.Item(7, 11).Value = 647  # was 645
$ws.Cells.Item(8, 11).Value = 1447  # was 1436
$ws.Cells.Item(10, 11).Value = 129  # was 128
$ws.Cells.Item(11, 11).Value = 409  # was 407
$ws.Cells.Item(15, 11).Value = 224  # was 221
$ws.Cells.Item(19, 11).Value = 641  # was 640
$ws.Cells.Item(20, 11).Value = 522  # was 521
$ws.Cells.Item(22, 11).Value = 67  # was 66
$ws.Cells.Item(27, 11).Value = 207  # was 208
$ws.Cells.Item(29, 11).Value = 1189  # was 1190
$ws.Cells.Item(30, 11).Value = 83  # was 81
$ws.Cells.Item(31, 11).Value = 245  # was 244
$ws.Cells.Item(33, 11).Value = 961  # was 959
$ws.Cells.Item(37, 11).Value = 748  # was 743
$ws.Cells.Item(41, 11).Value = 154  # was 153
$ws.Cells.Item(42, 11).Value = 814  # was 810
$ws.Cells.Item(43, 11).Value = 180  # was 179
$ws.Cells.Item(44, 11).Value = 182  # was 181
$ws.Cells.Item(46, 11).Value = 44  # was 43
$ws.Cells.Item(48, 11).Value = 275  # was 273
$ws.Cells.Item(51, 11).Value = 282  # was 281
$ws.Cells.Item(52, 11).Value = 581  # was 576
$ws.Cells.Item(54, 11).Value = 432  # was 431
$ws.Cells.Item(63, 11).Value = 59  # was 58
$ws.Cells.Item(67, 11).Value = 862  # was 858
$ws.Cells.Item(72, 11).Value = 113  # was 111
$ws.Cells.Item(76, 11).Value = 300  # was 299
$ws.Cells.Item(77, 11).Value = 151  # was 150
$ws.Cells.Item(78, 11).Value = 250  # was 248
$ws.Cells.Item(79, 11).Value = 554  # was 552
$ws.Cells.Item(85, 11).Value = 1022  # was 1019
$ws.Cells.Item(88, 11).Value = 236  # was 233
$ws.Cells.Item(89, 11).Value = 324  # was 323
$ws.Cells.Item(91, 11).Value = 256  # was 252
$ws.Cells.Item(95, 11).Value = 362  # was 360
$ws.Cells.Item(96, 11).Value = 232  # was 231
$ws.Cells.Item(97, 11).Value = 175  # was 172
$ws.Cells.Item(99, 11).Value = 362  # was 361
$ws.Cells.Item(101, 11).Value = 21997  # was 21915

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 11).Value = 81  # was 80
$ws.Cells.Item(7, 11).Value = 245  # was 244

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 11).Value = 316  # was 313
$ws.Cells.Item(6, 11).Value = 243  # was 242
$ws.Cells.Item(7, 11).Value = 862  # was 858

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 11).Value = 233  # was 232
$ws.Cells.Item(7, 11).Value = 432  # was 431

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(6, 11).Value = 342  # was 343
$ws.Cells.Item(7, 11).Value = 1189  # was 1190

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(4, 11).Value = 39  # was 38
$ws.Cells.Item(6, 11).Value = 129  # was 128
$ws.Cells.Item(7, 11).Value = 275  # was 273

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 11).Value = 193  # was 192
$ws.Cells.Item(7, 11).Value = 641  # was 640

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 11).Value = 51  # was 50
$ws.Cells.Item(7, 11).Value = 182  # was 181

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 11).Value = 19  # was 18
$ws.Cells.Item(7, 11).Value = 300  # was 299

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 11).Value = 43  # was 41
$ws.Cells.Item(7, 11).Value = 157  # was 155

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 11).Value = 52  # was 51
$ws.Cells.Item(7, 11).Value = 154  # was 153

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 11).Value = 221  # was 219
$ws.Cells.Item(3, 11).Value = 248  # was 246
$ws.Cells.Item(7, 11).Value = 814  # was 810

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 11).Value = 58  # was 57
$ws.Cells.Item(7, 11).Value = 129  # was 128

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 11).Value = 75  # was 74
$ws.Cells.Item(4, 11).Value = 23  # was 22
$ws.Cells.Item(7, 11).Value = 250  # was 248

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(6, 11).Value = 14  # was 13
$ws.Cells.Item(7, 11).Value = 44  # was 43

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 11).Value = 98  # was 97
$ws.Cells.Item(7, 11).Value = 232  # was 231

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 11).Value = 65  # was 63
$ws.Cells.Item(4, 11).Value = 13  # was 11
$ws.Cells.Item(7, 11).Value = 256  # was 252

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 11).Value = 179  # was 178
$ws.Cells.Item(6, 11).Value = 140  # was 139
$ws.Cells.Item(7, 11).Value = 554  # was 552

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 181  # was 180
$ws.Cells.Item(7, 11).Value = 522  # was 521

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 211  # was 210
$ws.Cells.Item(6, 11).Value = 178  # was 177
$ws.Cells.Item(7, 11).Value = 647  # was 645

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 11).Value = 82  # was 80
$ws.Cells.Item(6, 11).Value = 68  # was 67
$ws.Cells.Item(7, 11).Value = 224  # was 221

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(6, 11).Value = 134  # was 132
$ws.Cells.Item(7, 11).Value = 409  # was 407

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 11).Value = 61  # was 60
$ws.Cells.Item(7, 11).Value = 191  # was 190

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(2, 11).Value = 36  # was 35
$ws.Cells.Item(3, 11).Value = 37  # was 35
$ws.Cells.Item(7, 11).Value = 175  # was 172

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 11).Value = 61  # was 60
$ws.Cells.Item(3, 11).Value = 71  # was 70
$ws.Cells.Item(6, 11).Value = 96  # was 95
$ws.Cells.Item(7, 11).Value = 236  # was 233

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 11).Value = 99  # was 98
$ws.Cells.Item(7, 11).Value = 324  # was 323

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(2, 11).Value = 14  # was 13
$ws.Cells.Item(7, 11).Value = 57  # was 56

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 11).Value = 54  # was 55
$ws.Cells.Item(7, 11).Value = 207  # was 208

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 11).Value = 78  # was 77
$ws.Cells.Item(7, 11).Value = 282  # was 281

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(3, 11).Value = 47  # was 46
$ws.Cells.Item(7, 11).Value = 180  # was 179

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 335  # was 334
$ws.Cells.Item(3, 11).Value = 352  # was 350
$ws.Cells.Item(7, 11).Value = 1022  # was 1019

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(3, 11).Value = 20  # was 19
$ws.Cells.Item(7, 11).Value = 67  # was 66

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(2, 11).Value = 25  # was 24
$ws.Cells.Item(6, 11).Value = 54  # was 53
$ws.Cells.Item(7, 11).Value = 113  # was 111

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(4, 11).Value = 11  # was 10
$ws.Cells.Item(7, 11).Value = 151  # was 150

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 11).Value = 153  # was 151
$ws.Cells.Item(3, 11).Value = 166  # was 165
$ws.Cells.Item(4, 11).Value = 33  # was 32
$ws.Cells.Item(6, 11).Value = 210  # was 209
$ws.Cells.Item(7, 11).Value = 581  # was 576
